$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two tank description strings (G8 and G9) to reflect the new design:
# tanks now have a single powerful gun with medium rate of fire, instead of two turrets.
$ws.Range("G8").Value = "Большое количество здоровья. Имеет мощную пушку со средней скорострельностью. Поворотливый, но медленный."
$ws.Range("G9").Value = "Огромное количество здоровья. Имеет мощную пушку со средней скорострельностью. Поворотливый, но медленный."

# Update the view state: scroll the window so row 6 is at the top,
# and move the selection to G8
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("G8").Select()
